# Auto-generated edit script: update cryptocurrency price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '26.445.42'
$c.Style = "Normal"

$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '1.699.18'
$c.Style = "Normal"
$c = $ws.Range('E3')
$c.NumberFormat = "@"
$c.Value = '  +0.97%  '
$c.Style = "Normal"

$c = $ws.Range('E4')
$c.NumberFormat = "@"
$c.Value = '  +0.26%  '
$c.Style = "Normal"

$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '219.03'
$c.Style = "Normal"
$c = $ws.Range('E5')
$c.NumberFormat = "@"
$c.Value = '  +0.32%  '
$c.Style = "Normal"

$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '0.5475'
$c.Style = "Normal"
$c = $ws.Range('E6')
$c.NumberFormat = "@"
$c.Value = '  +4.01%  '
$c.Style = "Normal"

$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '1.010'
$c.Style = "Normal"
$c = $ws.Range('E7')
$c.NumberFormat = "@"
$c.Value = '  +0.22%  '
$c.Style = "Normal"

$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.2740'
$c.Style = "Normal"
$c = $ws.Range('E8')
$c.NumberFormat = "@"
$c.Value = '  +1.39%  '
$c.Style = "Normal"

$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.06447'
$c.Style = "Normal"
$c = $ws.Range('E9')
$c.NumberFormat = "@"
$c.Value = '  +0.42%  '
$c.Style = "Normal"

$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '22.05'
$c.Style = "Normal"
$c = $ws.Range('E10')
$c.NumberFormat = "@"
$c.Value = '  +0.04%  '
$c.Style = "Normal"

$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.07715'
$c.Style = "Normal"
$c = $ws.Range('E11')
$c.NumberFormat = "@"
$c.Value = '  +3.10%  '
$c.Style = "Normal"

$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '1.695.67'
$c.Style = "Normal"
$c = $ws.Range('E12')
$c.NumberFormat = "@"
$c.Value = '  +0.46%  '
$c.Style = "Normal"

$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '4.550'
$c.Style = "Normal"
$c = $ws.Range('E13')
$c.NumberFormat = "@"
$c.Value = '  +0.02%  '
$c.Style = "Normal"

$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '0.5836'
$c.Style = "Normal"
$c = $ws.Range('E14')
$c.NumberFormat = "@"
$c.Value = '  +0.42%  '
$c.Style = "Normal"

$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '0.000008408'
$c.Style = "Normal"
$c = $ws.Range('E15')
$c.NumberFormat = "@"
$c.Value = '  -1.00%  '
$c.Style = "Normal"

$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '65.79'
$c.Style = "Normal"
$c = $ws.Range('E16')
$c.NumberFormat = "@"
$c.Value = '  +2.32%  '
$c.Style = "Normal"

$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '26.466.58'
$c.Style = "Normal"
$c = $ws.Range('E17')
$c.NumberFormat = "@"
$c.Value = '  +0.50%  '
$c.Style = "Normal"

$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '4.948'
$c.Style = "Normal"
$c = $ws.Range('E18')
$c.NumberFormat = "@"
$c.Value = '  +0.52%  '
$c.Style = "Normal"

$c = $ws.Range('E19')
$c.NumberFormat = "@"
$c.Value = '  +0.30%  '
$c.Style = "Normal"

$c = $ws.Range('E20')
$c.NumberFormat = "@"
$c.Value = '  +1.10%  '
$c.Style = "Normal"

$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '191.61'
$c.Style = "Normal"
$c = $ws.Range('E21')
$c.NumberFormat = "@"
$c.Value = '  +1.42%  '
$c.Style = "Normal"

$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '6.261'
$c.Style = "Normal"
$c = $ws.Range('E22')
$c.NumberFormat = "@"
$c.Value = '  +1.10%  '
$c.Style = "Normal"

$c = $ws.Range('E23')
$c.NumberFormat = "@"
$c.Value = '  +0.31%  '
$c.Style = "Normal"

$c = $ws.Range('E24')
$c.NumberFormat = "@"
$c.Value = '  +3.38%  '
$c.Style = "Normal"

$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '0.1327'
$c.Style = "Normal"
$c = $ws.Range('E25')
$c.NumberFormat = "@"
$c.Value = '  +7.21%  '
$c.Style = "Normal"

$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '7.906'
$c.Style = "Normal"

$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '15.85'
$c.Style = "Normal"
$c = $ws.Range('E27')
$c.NumberFormat = "@"
$c.Value = '  +0.50%  '
$c.Style = "Normal"

$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '0.06249'
$c.Style = "Normal"
$c = $ws.Range('E28')
$c.NumberFormat = "@"
$c.Value = '  -5.73%  '
$c.Style = "Normal"

$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '1.379'
$c.Style = "Normal"
$c = $ws.Range('E29')
$c.NumberFormat = "@"
$c.Value = '  +2.45%  '
$c.Style = "Normal"

$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '1.333'
$c.Style = "Normal"
$c = $ws.Range('E30')
$c.NumberFormat = "@"
$c.Value = '  +0.33%  '
$c.Style = "Normal"

$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '3.606'
$c.Style = "Normal"
$c = $ws.Range('E31')
$c.NumberFormat = "@"
$c.Value = '  +0.86%  '
$c.Style = "Normal"

$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '3.614'
$c.Style = "Normal"
$c = $ws.Range('E32')
$c.NumberFormat = "@"
$c.Value = '  +1.53%  '
$c.Style = "Normal"

$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '1.694'
$c.Style = "Normal"
$c = $ws.Range('E33')
$c.NumberFormat = "@"
$c.Value = '  +1.75%  '
$c.Style = "Normal"

$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '1.041'
$c.Style = "Normal"
$c = $ws.Range('E34')
$c.NumberFormat = "@"
$c.Value = '  +1.28%  '
$c.Style = "Normal"

$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '0.6181'
$c.Style = "Normal"
$c = $ws.Range('E35')
$c.NumberFormat = "@"
$c.Value = '  -0.21%  '
$c.Style = "Normal"

$c = $ws.Range('E36')
$c.NumberFormat = "@"
$c.Value = '  +0.69%  '
$c.Style = "Normal"

$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '2.771'
$c.Style = "Normal"
$c = $ws.Range('E37')
$c.NumberFormat = "@"
$c.Value = '  +2.61%  '
$c.Style = "Normal"

$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '0.01642'
$c.Style = "Normal"

$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '1.118.18'
$c.Style = "Normal"
$c = $ws.Range('E39')
$c.NumberFormat = "@"
$c.Value = '  +0.76%  '
$c.Style = "Normal"

$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '6.147'
$c.Style = "Normal"
$c = $ws.Range('E40')
$c.NumberFormat = "@"
$c.Value = '  -4.06%  '
$c.Style = "Normal"

$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '0.8796'
$c.Style = "Normal"
$c = $ws.Range('E41')
$c.NumberFormat = "@"
$c.Value = '  +0.45%  '
$c.Style = "Normal"

$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '1.016'
$c.Style = "Normal"
$c = $ws.Range('E42')
$c.NumberFormat = "@"
$c.Value = '  +0.18%  '
$c.Style = "Normal"

$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '101.41'
$c.Style = "Normal"
$c = $ws.Range('E43')
$c.NumberFormat = "@"
$c.Value = '  +0.86%  '
$c.Style = "Normal"

$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '1.849.99'
$c.Style = "Normal"
$c = $ws.Range('E44')
$c.NumberFormat = "@"
$c.Value = '  +1.02%  '
$c.Style = "Normal"

$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '0.00000000109'
$c.Style = "Normal"
$c = $ws.Range('E45')
$c.NumberFormat = "@"
$c.Value = '  -4.43%  '
$c.Style = "Normal"

$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '57.56'
$c.Style = "Normal"
$c = $ws.Range('E46')
$c.NumberFormat = "@"
$c.Value = '  +1.36%  '
$c.Style = "Normal"

$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '8.232'
$c.Style = "Normal"
$c = $ws.Range('E47')
$c.NumberFormat = "@"
$c.Value = '  +0.93%  '
$c.Style = "Normal"

$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '1.007'
$c.Style = "Normal"
$c = $ws.Range('E48')
$c.NumberFormat = "@"
$c.Value = '  +0.33%  '
$c.Style = "Normal"

$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '0.05287'
$c.Style = "Normal"
$c = $ws.Range('E49')
$c.NumberFormat = "@"
$c.Value = '  +0.34%  '
$c.Style = "Normal"

$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '6.163'
$c.Style = "Normal"
$c = $ws.Range('E50')
$c.NumberFormat = "@"
$c.Value = '  +2.07%  '
$c.Style = "Normal"

$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '0.4303'
$c.Style = "Normal"
$c = $ws.Range('E51')
$c.NumberFormat = "@"
$c.Value = '  +0.01%  '
$c.Style = "Normal"

